$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Restructure header merges ---
$ws.Range("A1:D1").UnMerge()
$ws.Range("E1:H1").UnMerge()

# --- Step 2: Clear cells that no longer hold data ---
$ws.Range("A1").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("E1").ClearContents()
$ws.Range("E1").ClearFormats()
$ws.Range("E3").ClearContents()
$ws.Range("E3").ClearFormats()
$ws.Range("F3").ClearContents()
$ws.Range("F3").ClearFormats()
$ws.Range("G3").ClearContents()
$ws.Range("G3").ClearFormats()
$ws.Range("E4").ClearContents()
$ws.Range("E4").ClearFormats()
$ws.Range("F4").ClearContents()
$ws.Range("F4").ClearFormats()
$ws.Range("G4").ClearContents()
$ws.Range("G4").ClearFormats()
$ws.Range("E5").ClearContents()
$ws.Range("E5").ClearFormats()
$ws.Range("F5").ClearContents()
$ws.Range("F5").ClearFormats()
$ws.Range("G5").ClearContents()
$ws.Range("G5").ClearFormats()
$ws.Range("E6").ClearContents()
$ws.Range("E6").ClearFormats()
$ws.Range("F6").ClearContents()
$ws.Range("F6").ClearFormats()
$ws.Range("G6").ClearContents()
$ws.Range("G6").ClearFormats()
$ws.Range("E7").ClearContents()
$ws.Range("E7").ClearFormats()
$ws.Range("F7").ClearContents()
$ws.Range("F7").ClearFormats()
$ws.Range("G7").ClearContents()
$ws.Range("G7").ClearFormats()
$ws.Range("E8").ClearContents()
$ws.Range("E8").ClearFormats()
$ws.Range("F8").ClearContents()
$ws.Range("F8").ClearFormats()
$ws.Range("G8").ClearContents()
$ws.Range("G8").ClearFormats()
$ws.Range("E9").ClearContents()
$ws.Range("E9").ClearFormats()
$ws.Range("F9").ClearContents()
$ws.Range("F9").ClearFormats()
$ws.Range("G9").ClearContents()
$ws.Range("G9").ClearFormats()
$ws.Range("E10").ClearContents()
$ws.Range("E10").ClearFormats()
$ws.Range("F10").ClearContents()
$ws.Range("F10").ClearFormats()
$ws.Range("G10").ClearContents()
$ws.Range("G10").ClearFormats()
$ws.Range("E11").ClearContents()
$ws.Range("E11").ClearFormats()
$ws.Range("F11").ClearContents()
$ws.Range("F11").ClearFormats()
$ws.Range("G11").ClearContents()
$ws.Range("G11").ClearFormats()
$ws.Range("E12").ClearContents()
$ws.Range("E12").ClearFormats()
$ws.Range("F12").ClearContents()
$ws.Range("F12").ClearFormats()
$ws.Range("G12").ClearContents()
$ws.Range("G12").ClearFormats()
$ws.Range("E13").ClearContents()
$ws.Range("E13").ClearFormats()
$ws.Range("F13").ClearContents()
$ws.Range("F13").ClearFormats()
$ws.Range("G13").ClearContents()
$ws.Range("G13").ClearFormats()
$ws.Range("E14").ClearContents()
$ws.Range("E14").ClearFormats()
$ws.Range("F14").ClearContents()
$ws.Range("F14").ClearFormats()
$ws.Range("G14").ClearContents()
$ws.Range("G14").ClearFormats()

# --- Step 3: Set header date cells (text, bold + centered) ---
$ws.Range("B1").Value2 = "'05-07-2023"
$ws.Range("B1").ClearFormats()
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("F1").Value2 = "'05-07-2023"
$ws.Range("F1").ClearFormats()
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("J1").Value2 = "'05-07-2023"
$ws.Range("J1").ClearFormats()
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108

# --- Step 4: Re-merge header ranges ---
$ws.Range("B1:E1").Merge()
$ws.Range("F1:I1").Merge()
$ws.Range("J1:M1").Merge()

# --- Step 5: Populate data rows 2-14 ---
$ws.Range("B2").Value2 = "'+14.13%"
$ws.Range("B2").ClearFormats()

$ws.Range("D2").ClearFormats()
$ws.Range("D2").Value2 = 9.289999999999999
$ws.Range("D2").HorizontalAlignment = -4108

$ws.Range("E2").ClearFormats()
$ws.Range("E2").Value2 = ""
$ws.Range("E2").HorizontalAlignment = -4108

$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value2 = ""
$ws.Range("F2").HorizontalAlignment = -4108

$ws.Range("G2").ClearFormats()
$ws.Range("G2").Value2 = ""
$ws.Range("G2").HorizontalAlignment = -4108

$ws.Range("I2").Value2 = "'1er Janvier"
$ws.Range("I2").ClearFormats()
$ws.Range("I2").HorizontalAlignment = -4108

$ws.Range("J2").Value2 = "'+14.13%"
$ws.Range("J2").ClearFormats()
$ws.Range("J2").HorizontalAlignment = -4108

$ws.Range("K2").ClearFormats()
$ws.Range("K2").Value2 = 11.87
$ws.Range("K2").HorizontalAlignment = -4108

$ws.Range("L2").ClearFormats()
$ws.Range("L2").Value2 = 9.289999999999999
$ws.Range("L2").HorizontalAlignment = -4108

$ws.Range("M2").Value2 = "'1er Janvier"
$ws.Range("M2").ClearFormats()
$ws.Range("M2").HorizontalAlignment = -4108

$ws.Range("N2").Value2 = "'+14.13%"
$ws.Range("N2").ClearFormats()
$ws.Range("N2").HorizontalAlignment = -4108

$ws.Range("O2").ClearFormats()
$ws.Range("O2").Value2 = 11.87

$ws.Range("P2").ClearFormats()
$ws.Range("P2").Value2 = 9.289999999999999

$ws.Range("B3").Value2 = "'-1.51%"
$ws.Range("B3").ClearFormats()

$ws.Range("D3").ClearFormats()
$ws.Range("D3").Value2 = 10.51

$ws.Range("I3").Value2 = "'1 semaine"
$ws.Range("I3").ClearFormats()

$ws.Range("J3").Value2 = "'-1.51%"
$ws.Range("J3").ClearFormats()

$ws.Range("K3").ClearFormats()
$ws.Range("K3").Value2 = 10.81

$ws.Range("L3").ClearFormats()
$ws.Range("L3").Value2 = 10.51

$ws.Range("M3").Value2 = "'1 semaine"
$ws.Range("M3").ClearFormats()

$ws.Range("N3").Value2 = "'-1.51%"
$ws.Range("N3").ClearFormats()

$ws.Range("O3").ClearFormats()
$ws.Range("O3").Value2 = 10.81

$ws.Range("P3").ClearFormats()
$ws.Range("P3").Value2 = 10.51

$ws.Range("B4").Value2 = "'-0.88%"
$ws.Range("B4").ClearFormats()

$ws.Range("C4").ClearFormats()
$ws.Range("C4").Value2 = 10.81

$ws.Range("D4").ClearFormats()
$ws.Range("D4").Value2 = 10.21

$ws.Range("I4").Value2 = "'1 mois"
$ws.Range("I4").ClearFormats()

$ws.Range("J4").Value2 = "'-0.88%"
$ws.Range("J4").ClearFormats()

$ws.Range("K4").ClearFormats()
$ws.Range("K4").Value2 = 10.81

$ws.Range("L4").ClearFormats()
$ws.Range("L4").Value2 = 10.21

$ws.Range("M4").Value2 = "'1 mois"
$ws.Range("M4").ClearFormats()

$ws.Range("N4").Value2 = "'-0.88%"
$ws.Range("N4").ClearFormats()

$ws.Range("O4").ClearFormats()
$ws.Range("O4").Value2 = 10.81

$ws.Range("P4").ClearFormats()
$ws.Range("P4").Value2 = 10.21

$ws.Range("B5").Value2 = "'-7.57%"
$ws.Range("B5").ClearFormats()

$ws.Range("D5").ClearFormats()
$ws.Range("D5").Value2 = 10.21

$ws.Range("I5").Value2 = "'3 mois"
$ws.Range("I5").ClearFormats()

$ws.Range("J5").Value2 = "'-7.57%"
$ws.Range("J5").ClearFormats()

$ws.Range("K5").ClearFormats()
$ws.Range("K5").Value2 = 11.87

$ws.Range("L5").ClearFormats()
$ws.Range("L5").Value2 = 10.21

$ws.Range("M5").Value2 = "'3 mois"
$ws.Range("M5").ClearFormats()

$ws.Range("N5").Value2 = "'-7.57%"
$ws.Range("N5").ClearFormats()

$ws.Range("O5").ClearFormats()
$ws.Range("O5").Value2 = 11.87

$ws.Range("P5").ClearFormats()
$ws.Range("P5").Value2 = 10.21

$ws.Range("B6").Value2 = "'+7.09%"
$ws.Range("B6").ClearFormats()

$ws.Range("D6").ClearFormats()
$ws.Range("D6").Value2 = 9.49

$ws.Range("I6").Value2 = "'6 mois"
$ws.Range("I6").ClearFormats()

$ws.Range("J6").Value2 = "'+7.09%"
$ws.Range("J6").ClearFormats()

$ws.Range("K6").ClearFormats()
$ws.Range("K6").Value2 = 11.87

$ws.Range("L6").ClearFormats()
$ws.Range("L6").Value2 = 9.49

$ws.Range("M6").Value2 = "'6 mois"
$ws.Range("M6").ClearFormats()

$ws.Range("N6").Value2 = "'+7.09%"
$ws.Range("N6").ClearFormats()

$ws.Range("O6").ClearFormats()
$ws.Range("O6").Value2 = 11.87

$ws.Range("P6").ClearFormats()
$ws.Range("P6").Value2 = 9.49

$ws.Range("B7").Value2 = "'-3.83%"
$ws.Range("B7").ClearFormats()

$ws.Range("D7").ClearFormats()
$ws.Range("D7").Value2 = 9.08

$ws.Range("I7").Value2 = "'1 an"
$ws.Range("I7").ClearFormats()

$ws.Range("J7").Value2 = "'-3.83%"
$ws.Range("J7").ClearFormats()

$ws.Range("K7").ClearFormats()
$ws.Range("K7").Value2 = 11.87

$ws.Range("L7").ClearFormats()
$ws.Range("L7").Value2 = 9.08

$ws.Range("M7").Value2 = "'1 an"
$ws.Range("M7").ClearFormats()

$ws.Range("N7").Value2 = "'-3.83%"
$ws.Range("N7").ClearFormats()

$ws.Range("O7").ClearFormats()
$ws.Range("O7").Value2 = 11.87

$ws.Range("P7").ClearFormats()
$ws.Range("P7").Value2 = 9.08

$ws.Range("B8").Value2 = "'-3.27%"
$ws.Range("B8").ClearFormats()

$ws.Range("D8").ClearFormats()
$ws.Range("D8").Value2 = 8.630000000000001

$ws.Range("I8").Value2 = "'3 ans"
$ws.Range("I8").ClearFormats()

$ws.Range("J8").Value2 = "'-3.27%"
$ws.Range("J8").ClearFormats()

$ws.Range("K8").ClearFormats()
$ws.Range("K8").Value2 = 11.94

$ws.Range("L8").ClearFormats()
$ws.Range("L8").Value2 = 8.630000000000001

$ws.Range("M8").Value2 = "'3 ans"
$ws.Range("M8").ClearFormats()

$ws.Range("N8").Value2 = "'-3.27%"
$ws.Range("N8").ClearFormats()

$ws.Range("O8").ClearFormats()
$ws.Range("O8").Value2 = 11.94

$ws.Range("P8").ClearFormats()
$ws.Range("P8").Value2 = 8.630000000000001

$ws.Range("B9").Value2 = "'-28.72%"
$ws.Range("B9").ClearFormats()

$ws.Range("D9").ClearFormats()
$ws.Range("D9").Value2 = 8.630000000000001

$ws.Range("I9").Value2 = "'5 ans"
$ws.Range("I9").ClearFormats()

$ws.Range("J9").Value2 = "'-28.72%"
$ws.Range("J9").ClearFormats()

$ws.Range("K9").ClearFormats()
$ws.Range("K9").Value2 = 15.38

$ws.Range("L9").ClearFormats()
$ws.Range("L9").Value2 = 8.630000000000001

$ws.Range("M9").Value2 = "'5 ans"
$ws.Range("M9").ClearFormats()

$ws.Range("N9").Value2 = "'-28.72%"
$ws.Range("N9").ClearFormats()

$ws.Range("O9").ClearFormats()
$ws.Range("O9").Value2 = 15.38

$ws.Range("P9").ClearFormats()
$ws.Range("P9").Value2 = 8.630000000000001

$ws.Range("B10").Value2 = "'+46.81%"
$ws.Range("B10").ClearFormats()

$ws.Range("D10").ClearFormats()
$ws.Range("D10").Value2 = 7.1

$ws.Range("I10").Value2 = "'10 ans"
$ws.Range("I10").ClearFormats()

$ws.Range("J10").Value2 = "'+46.81%"
$ws.Range("J10").ClearFormats()

$ws.Range("K10").ClearFormats()
$ws.Range("K10").Value2 = 16.98

$ws.Range("L10").ClearFormats()
$ws.Range("L10").Value2 = 7.1

$ws.Range("M10").Value2 = "'10 ans"
$ws.Range("M10").ClearFormats()

$ws.Range("N10").Value2 = "'+46.81%"
$ws.Range("N10").ClearFormats()

$ws.Range("O10").ClearFormats()
$ws.Range("O10").Value2 = 16.98

$ws.Range("P10").ClearFormats()
$ws.Range("P10").Value2 = 7.1

$ws.Range("D11").ClearFormats()
$ws.Range("D11").Value2 = 10.46

$ws.Range("I11").Value2 = "'MM20"
$ws.Range("I11").ClearFormats()

$ws.Range("L11").ClearFormats()
$ws.Range("L11").Value2 = 10.46

$ws.Range("M11").Value2 = "'MM20"
$ws.Range("M11").ClearFormats()

$ws.Range("N11").ClearFormats()
$ws.Range("N11").Value2 = ""

$ws.Range("O11").ClearFormats()
$ws.Range("O11").Value2 = ""

$ws.Range("P11").ClearFormats()
$ws.Range("P11").Value2 = 10.46

$ws.Range("D12").ClearFormats()
$ws.Range("D12").Value2 = 11.02

$ws.Range("I12").Value2 = "'MM50"
$ws.Range("I12").ClearFormats()

$ws.Range("L12").ClearFormats()
$ws.Range("L12").Value2 = 11.02

$ws.Range("M12").Value2 = "'MM50"
$ws.Range("M12").ClearFormats()

$ws.Range("N12").ClearFormats()
$ws.Range("N12").Value2 = ""

$ws.Range("O12").ClearFormats()
$ws.Range("O12").Value2 = ""

$ws.Range("P12").ClearFormats()
$ws.Range("P12").Value2 = 11.02

$ws.Range("D13").ClearFormats()
$ws.Range("D13").Value2 = 10.98

$ws.Range("I13").Value2 = "'MM100"
$ws.Range("I13").ClearFormats()

$ws.Range("L13").ClearFormats()
$ws.Range("L13").Value2 = 10.98

$ws.Range("M13").Value2 = "'MM100"
$ws.Range("M13").ClearFormats()

$ws.Range("N13").ClearFormats()
$ws.Range("N13").Value2 = ""

$ws.Range("O13").ClearFormats()
$ws.Range("O13").Value2 = ""

$ws.Range("P13").ClearFormats()
$ws.Range("P13").Value2 = 10.98

$ws.Range("D14").ClearFormats()
$ws.Range("D14").Value2 = 65.53

$ws.Range("I14").Value2 = "'RSI14"
$ws.Range("I14").ClearFormats()

$ws.Range("L14").ClearFormats()
$ws.Range("L14").Value2 = 65.53

$ws.Range("M14").Value2 = "'RSI14"
$ws.Range("M14").ClearFormats()

$ws.Range("N14").ClearFormats()
$ws.Range("N14").Value2 = ""

$ws.Range("O14").ClearFormats()
$ws.Range("O14").Value2 = ""

$ws.Range("P14").ClearFormats()
$ws.Range("P14").Value2 = 65.53

# --- Step 6: Force row 15 into used range (matches original empty trailing row) ---
$ws.Range("A15").Borders.LineStyle = 0
